$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D holds the "Speaker" values (header in D1).
# Rename speaker tags: "RBD" -> "T" and "Student" -> "S" throughout the used range.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -eq "RBD") {
        $cell.Value = "T"
    } elseif ($val -eq "Student") {
        $cell.Value = "S"
    }
}
